$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.373.52"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "'1.868.32"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'234.67"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4700"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").Value = "'0.2873"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "'21.53"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'0.07877"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "'96.95"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "'1.871.70"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'0.6917"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "'5.109"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "'268.41"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'30.345.23"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "'13.95"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").Value = "'0.000007691"
$ws.Range("E19").Value = "  +4.00%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'2.116.54"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'5.246"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'6.197"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'9.401"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").Value = "'167.64"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'18.87"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'1.950"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.09938"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.361"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "'4.410"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "'1.462"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "'4.067"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "'0.04753"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'0.7054"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "'2.719"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "'2.798"
$ws.Range("E39").Value = "  +6.81%  "
$ws.Range("D40").Value = "'6.280"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "'73.96"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "'1.952"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8437"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4173"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'103.33"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").Value = "'987.96"
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("D48").Value = "'7.116"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "'9.121"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").Value = "'34.56"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").Value = "'0.05679"
$ws.Range("E51").Value = "  +0.30%  "
